$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''38.746.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.49%  '

$ws.Range("D3").Value = '''2.098.41'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.27%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '''228.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.07%  '

$ws.Range("D6").Value = '''0.618'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.84%  '

$ws.Range("D7").Value = '''61.50'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.83%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").Value = '  +1.90%  '

$ws.Range("D10").Value = '''0.0843'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.66%  '

$ws.Range("D11").Value = '''0.104'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.03%  '

$ws.Range("D12").Value = '''15.51'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.60%  '

$ws.Range("D13").Value = '''2.411.27'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.36%  '

$ws.Range("D14").Value = '''21.98'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.22%  '

$ws.Range("D15").Value = '''0.807'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.86%  '

$ws.Range("D16").Value = '''5.50'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.19%  '

$ws.Range("D17").Value = '''2.110.44'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.63%  '

$ws.Range("D18").Value = '''38.748.38'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.70%  '

$ws.Range("D19").Value = '''71.97'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.57%  '

$ws.Range("D20").Value = '''6.08'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.03%  '

$ws.Range("D21").Value = '''0.0₃0839'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.33%  '

$ws.Range("D22").Value = '''227.53'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.57%  '

$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("E24").Value = '  -2.20%  '

$ws.Range("D26").Value = '''171.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.06%  '

$ws.Range("D27").Value = '''9.55'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.96%  '

$ws.Range("D28").Value = '''0.137'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.53%  '

$ws.Range("E29").Value = '  +3.67%  '

$ws.Range("D30").Value = '''19.32'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.90%  '

$ws.Range("E31").Value = '  +3.58%  '

$ws.Range("E32").Value = '  +1.02%  '

$ws.Range("D33").Value = '''4.53'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.21%  '

$ws.Range("D34").Value = '''4.76'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.38%  '

$ws.Range("E35").Value = '  +2.43%  '

$ws.Range("E36").Value = '  +2.24%  '

$ws.Range("D37").Value = '''2.40'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.42%  '

$ws.Range("E38").Value = '  +1.51%  '

$ws.Range("E39").Value = '  -0.07%  '

$ws.Range("D40").Value = '''18.20'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.89%  '

$ws.Range("E41").Value = '  +4.34%  '

$ws.Range("D42").Value = '''101.40'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.30%  '

$ws.Range("D43").Value = '''1.533.48'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.41%  '

$ws.Range("E44").Value = '  -0.90%  '

$ws.Range("E45").Value = '  +4.36%  '

$ws.Range("D46").Value = '''0.0909'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.49%  '

$ws.Range("E47").Value = '  +1.96%  '

$ws.Range("D48").Value = '''4.12'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.04%  '

$ws.Range("E49").Value = '  +1.63%  '

$ws.Range("E50").Value = '  -1.04%  '

$ws.Range("D51").Value = '''2.294.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.20%  '
